$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-6, columns A-I (J stays the shared string "train_dim2_1")
$data = @(
    @(1, 4, 2, 7, 5, 3, 3, 43, 5),
    @(2, 0, 0, 5, 1, 5, 1, 65, 5),
    @(3, 3, 1, 4, 6, 1, 5, 21, 5),
    @(4, 0, 2, 4, 4, 4, 2, 54, 5),
    @(5, 4, 1, 6, 5, 2, 4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$wb.Save()
